$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1487266666666667
$ws.Range("H2").Value = 0.44618
$ws.Range("I2").Value = 0.1109321277273384
$ws.Range("J2").Value = 0.1109321277273384
$ws.Range("M2").Value = 0.242595
$ws.Range("N2").Value = 0.727785
$ws.Range("O2").Value = 0.03230700759563258
$ws.Range("P2").Value = 0.03230700759563257
$ws.Range("Q2").Value = 0.0360803457
$ws.Range("R2").Value = 0.3247231113
$ws.Range("S2").Value = 0.003583885093086805
$ws.Range("T2").Value = 0.003583885093086804
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1487266666666667
$ws.Range("H3").Value = 0.44618
$ws.Range("I3").Value = 0.1109321277273384
$ws.Range("J3").Value = 0.1109321277273384
$ws.Range("O3").Value = 0.4893229309549773
$ws.Range("P3").Value = 0.4893229309549771
$ws.Range("Q3").Value = 0.5464740259688889
$ws.Range("R3").Value = 4.918266233720001
$ws.Range("S3").Value = 0.05428163387661311
$ws.Range("T3").Value = 0.0542816338766131
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1487266666666667
$ws.Range("H4").Value = 0.44618
$ws.Range("I4").Value = 0.1109321277273384
$ws.Range("J4").Value = 0.1109321277273384
$ws.Range("M4").Value = 3.574634666666667
$ws.Range("N4").Value = 10.723904
$ws.Range("O4").Value = 0.4760434029044768
$ws.Range("P4").Value = 0.4760434029044767
$ws.Range("Q4").Value = 0.5316434985244445
$ws.Range("R4").Value = 4.784791486720001
$ws.Range("S4").Value = 0.05280850757475622
$ws.Range("T4").Value = 0.05280850757475622
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1487266666666667
$ws.Range("H5").Value = 0.44618
$ws.Range("I5").Value = 0.1109321277273384
$ws.Range("J5").Value = 0.1109321277273384
$ws.Range("M5").Value = 0.017471
$ws.Range("N5").Value = 0.052413
$ws.Range("O5").Value = 0.002326658544913526
$ws.Range("P5").Value = 0.002326658544913525
$ws.Range("Q5").Value = 0.002598403593333333
$ws.Range("R5").Value = 0.02338563234
$ws.Range("S5").Value = 0.0002581011828822505
$ws.Range("T5").Value = 0.0002581011828822505
$ws.Range("I6").Value = 0.05522216136400421
$ws.Range("J6").Value = 0.05522216136400422
$ws.Range("M6").Value = 0.242595
$ws.Range("N6").Value = 0.727785
$ws.Range("O6").Value = 0.03230700759563258
$ws.Range("P6").Value = 0.03230700759563257
$ws.Range("Q6").Value = 0.017960844285
$ws.Range("R6").Value = 0.161647598565
$ws.Range("S6").Value = 0.001784062786634132
$ws.Range("T6").Value = 0.001784062786634132
$ws.Range("I7").Value = 0.05522216136400421
$ws.Range("J7").Value = 0.05522216136400422
$ws.Range("O7").Value = 0.4893229309549773
$ws.Range("P7").Value = 0.4893229309549771
$ws.Range("S7").Value = 0.02702146985230324
$ws.Range("T7").Value = 0.02702146985230324
$ws.Range("I8").Value = 0.05522216136400421
$ws.Range("J8").Value = 0.05522216136400422
$ws.Range("M8").Value = 3.574634666666667
$ws.Range("N8").Value = 10.723904
$ws.Range("O8").Value = 0.4760434029044768
$ws.Range("P8").Value = 0.4760434029044767
$ws.Range("Q8").Value = 0.2646528437262222
$ws.Range("R8").Value = 2.381875593536
$ws.Range("S8").Value = 0.02628814561146069
$ws.Range("T8").Value = 0.02628814561146069
$ws.Range("I9").Value = 0.05522216136400421
$ws.Range("J9").Value = 0.05522216136400422
$ws.Range("M9").Value = 0.017471
$ws.Range("N9").Value = 0.052413
$ws.Range("O9").Value = 0.002326658544913526
$ws.Range("P9").Value = 0.002326658544913525
$ws.Range("Q9").Value = 0.001293488779666667
$ws.Range("R9").Value = 0.011641399017
$ws.Range("S9").Value = 0.000128483113606154
$ws.Range("T9").Value = 0.000128483113606154
$ws.Range("G10").Value = 1.117936666666667
$ws.Range("H10").Value = 3.35381
$ws.Range("I10").Value = 0.8338457109086573
$ws.Range("J10").Value = 0.8338457109086574
$ws.Range("M10").Value = 0.242595
$ws.Range("N10").Value = 0.727785
$ws.Range("O10").Value = 0.03230700759563258
$ws.Range("P10").Value = 0.03230700759563257
$ws.Range("Q10").Value = 0.2712058456500001
$ws.Range("R10").Value = 2.44085261085
$ws.Range("S10").Value = 0.02693905971591164
$ws.Range("T10").Value = 0.02693905971591164
$ws.Range("G11").Value = 1.117936666666667
$ws.Range("H11").Value = 3.35381
$ws.Range("I11").Value = 0.8338457109086573
$ws.Range("J11").Value = 0.8338457109086574
$ws.Range("O11").Value = 0.4893229309549773
$ws.Range("P11").Value = 0.4893229309549771
$ws.Range("Q11").Value = 4.10769208174889
$ws.Range("R11").Value = 36.96922873574
$ws.Range("S11").Value = 0.4080198272260608
$ws.Range("T11").Value = 0.4080198272260608
$ws.Range("G12").Value = 1.117936666666667
$ws.Range("H12").Value = 3.35381
$ws.Range("I12").Value = 0.8338457109086573
$ws.Range("J12").Value = 0.8338457109086574
$ws.Range("M12").Value = 3.574634666666667
$ws.Range("N12").Value = 10.723904
$ws.Range("O12").Value = 0.4760434029044768
$ws.Range("P12").Value = 0.4760434029044767
$ws.Range("Q12").Value = 3.996215163804445
$ws.Range("R12").Value = 35.96593647424
$ws.Range("S12").Value = 0.3969467497182598
$ws.Range("T12").Value = 0.3969467497182598
$ws.Range("G13").Value = 1.117936666666667
$ws.Range("H13").Value = 3.35381
$ws.Range("I13").Value = 0.8338457109086573
$ws.Range("J13").Value = 0.8338457109086574
$ws.Range("M13").Value = 0.017471
$ws.Range("N13").Value = 0.052413
$ws.Range("O13").Value = 0.002326658544913526
$ws.Range("P13").Value = 0.002326658544913525
$ws.Range("Q13").Value = 0.01953147150333333
$ws.Range("R13").Value = 0.17578324353
$ws.Range("S13").Value = 0.001940074248425121
$ws.Range("T13").Value = 0.001940074248425121
